# Finalized definition of virtual environment
# Update calibration figures across the start_price, Linear and NonLinear sheets.

$wb = $excel.ActiveWorkbook

# --- start_price sheet ---
$wsStart = $wb.Worksheets.Item("start_price")
$wsStart.Range("A2").Value = 13036.3733693968

# --- Linear sheet ---
$wsLinear = $wb.Worksheets.Item("Linear")
$wsLinear.Range("B2").Value = -0.8471131810676074
$wsLinear.Range("B3").Value = 0.1313170304116882
$wsLinear.Range("B4").Value = 31135.10500681926
$wsLinear.Range("B5").Value = "[1.0, 0.16666484744293045, 0.026888590862929895, 0.011586570939146721, 0.013439068188339546, -0.04388511630743954, 0.10310668537439181, 0.22375879637338877, 0.07125787748016232, -0.041473969659244964, -0.031453325315431914, -0.03366932283278228, -0.06763654787599654, 0.09154620980393374, 0.2013496382853902, 0.04136405623942855, -0.05135151979186593, -0.027021277824181324, -0.015006708274903698, -0.04329592805606688]"

# --- NonLinear sheet ---
$wsNonLinear = $wb.Worksheets.Item("NonLinear")
$wsNonLinear.Range("B3").Value = 0.985854189336235
$wsNonLinear.Range("B4").Value = -15.39084058107468
$wsNonLinear.Range("B5").Value = -0.06167806106399603
$wsNonLinear.Range("B6").Value = 28767.99384543934
$wsNonLinear.Range("B7").Value = -3.51249155586502
$wsNonLinear.Range("B8").Value = 0.2141775960811582
$wsNonLinear.Range("B9").Value = 33375.61582050646
$wsNonLinear.Range("B10").Value = "[0.9999999999999998, 0.16683640515536963, 0.029010434742807067, 0.010778014251443384, 0.007609802789236635, -0.04383622131140616, 0.10067858761429775, 0.21987930309605272, 0.07271597068785325, -0.043964385122186526, -0.033929523200823826, -0.036308186342196416, -0.06704920927135136, 0.0890182628218022, 0.20014405908826363, 0.04096481749374034, -0.05156274633814644, -0.03025149312566355, -0.01516062450994968, -0.04222097851607758]"
